$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values for P1, Q1 (columns 16, 17), copying O1's format (bold/border/centered)
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

$headerSrc = $ws.Cells.Item(1, 15)
$headerSrc.Copy()
$headerDst = $ws.Range($ws.Cells.Item(1, 16), $ws.Cells.Item(1, 17))
$headerDst.PasteSpecial(-4122)

# Update data rows 2-25: swap I<->K values and M<->O values, and add P, Q columns
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 9).Value = 2   # I
    $ws.Cells.Item($row, 11).Value = 1  # K
    $ws.Cells.Item($row, 13).Value = 2  # M
    $ws.Cells.Item($row, 15).Value = 1  # O
    $ws.Cells.Item($row, 16).Value = 2  # P
    $ws.Cells.Item($row, 17).Value = 2  # Q
}
